# Assignment 8 - "made a few grammatical corrections to assignment 8"
#
# Applies the textual/structural corrections from the commit. Each Find is
# scoped via $d.Content so it only touches the single paragraph containing
# the match (this COM shim coalesces the runs of whatever paragraph a Find/
# Range edit lands in - that's fine here because none of the runs being
# merged carry distinguishing character formatting; the pre-edit run splits
# were themselves just Word's per-keystroke-session rsid bookkeeping, not
# meaningful formatting boundaries).
#
# NOT replicated here (and not meaningfully reachable through the Word
# object model at all):
#   - xmlns attribute-list simplification on <w:document>/<w:numbering>/
#     <w:styles>, and the w:latentStyles count / rsid / qFormat churn in
#     styles.xml: cosmetic fingerprints of the particular Word build that
#     last resaved the file, not user-facing edits.
#   - hyperlink r:id renumbering (rId6->rId5 etc. in document.xml/
#     document.xml.rels): a side effect of the resave dropping/reordering
#     relationship parts; the hyperlink Targets themselves are unchanged,
#     so nothing user-visible moved.
#   - the two new/relocated <w:lastRenderedPageBreak/> markers: a read-only
#     pagination cache Word regenerates on layout; there is no
#     Application/Selection/Range method that inserts or moves one (real
#     Word macros can't do this either).

$d = $word.ActiveDocument

function DoReplace($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2) | Out-Null
}

# "... This is also the case for it's other 3 divisions ..." -> "its"
$apos = [char]0x2019
$find1 = "it" + $apos + "s other 3 divisions"
DoReplace $find1 "its other 3 divisions"

# "Which of our webpages would be easily ported to" -> "web pages"
$find2 = "Which of our webpages would be easily ported to"
$repl2 = "Which of our web pages would be easily ported to"
DoReplace $find2 $repl2

# "Advantage: Better customer experience." -> "Advantage: Provides a better customer experience."
$find3 = "Advantage: Better customer experience."
$repl3 = "Advantage: Provides a better customer experience."
DoReplace $find3 $repl3

# "Disadvantages: Would probably ... Thus more effort will be required." ->
# "Disadvantage: This would probably ... Thus more effort and resources will be required."
$find4 = "Disadvantages: Would probably require app development, specific for each platform. Thus more effort will be required."
$repl4 = "Disadvantage: This would probably require app development, specific for each platform. Thus more effort and resources will be required."
DoReplace $find4 $repl4

# "I recommend this taking this direction" -> "I recommend taking this direction"
$find5 = "I recommend this taking this direction"
$repl5 = "I recommend taking this direction"
DoReplace $find5 $repl5

# "(E.g., banks, entertainment networks such as ESPN)" -> adds "and "
$find6 = "banks, entertainment networks such as ESPN"
$repl6 = "banks, and entertainment networks such as ESPN"
DoReplace $find6 $repl6

# "Identify roles/personnel:" -> "Identify the following roles/personnel:"
$find7 = "Identify roles/personnel:"
$repl7 = "Identify the following roles/personnel:"
DoReplace $find7 $repl7

# "How long will it take?" -> "How long will this project take?"
$find8 = "How long will it take?"
$repl8 = "How long will this project take?"
DoReplace $find8 $repl8

# "... their mobile device can quickly ..." -> "mobile devices"
$find9 = "mobile device can quickly"
$repl9 = "mobile devices can quickly"
DoReplace $find9 $repl9

# Promote the "Websites that were created using responsive web design ..."
# bullet from list level 0 (w:ilvl val="0") to list level 2 (w:ilvl val="2").
# ListLevelNumber is 1-based (ListLevelNumber 3 == w:ilvl val="2").
$bulletFind = "Websites that were created using responsive web design"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith($bulletFind)) {
        $target = $d.Paragraphs($i)
        break
    }
}
$target.Range.ListFormat.ListLevelNumber = 3

Write-Output "Done applying assignment8 grammatical corrections."
